$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # 数据
$ws2 = $wb.Worksheets.Item(2)   # 说明

# ---------------------------------------------------------------------------
# Sheet "数据" (sheet1): fill in the newly-reported data for rows 16 and 17
# ---------------------------------------------------------------------------
$ws1.Range("C16").Value = 70
$ws1.Range("D16").Value = 6
$ws1.Range("E16").Value = 3
$ws1.Range("F16").Value = 1
$ws1.Range("G16").Value = 495
$ws1.Range("H16").Value = 23
$ws1.Range("I16").Value = 31
$ws1.Range("J16").Value = 441
$ws1.Range("K16").Value = 61
$ws1.Range("L16").Value = 9

$ws1.Range("C17").Value = 77
$ws1.Range("D17").Value = 15
$ws1.Range("E17").Value = 1
$ws1.Range("F17").Value = 7
$ws1.Range("G17").Value = 572
$ws1.Range("H17").Value = 38
$ws1.Range("I17").Value = 32
$ws1.Range("J17").Value = 502
$ws1.Range("K17").Value = 61
$ws1.Range("L17").Value = 16

# ---------------------------------------------------------------------------
# Sheet "说明" (sheet2): insert a new row documenting the lockdown announcement
# ---------------------------------------------------------------------------
$ws2.Rows.Item(3).Insert()
$ws2.Range("A3").Value = 43853
$ws2.Range("A3").Style = "Hyperlink"
$ws2.Range("A3").NumberFormat = "[$-411]ggge""年""m""月""d""日"""
$ws2.Range("B3").Value = "武汉市宣部封城"

# ---------------------------------------------------------------------------
# Active sheet / selections (last Select()/Activate() wins -> leave "数据" active)
# ---------------------------------------------------------------------------
$ws2.Range("D12").Select()
$ws1.Activate()
$ws1.Range("I22").Select()
